$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Reln"
$ws.Cells.Item(2,3).Value = "Vldlr"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0.03565433333333334
$ws.Cells.Item(2,8).Value = 0.106963
$ws.Cells.Item(2,9).Value = 0.002412342638581826
$ws.Cells.Item(2,10).Value = 0.002412342638581825
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.06089466666666667
$ws.Cells.Item(2,14).Value = 0.182684
$ws.Cells.Item(2,15).Value = 0.001903591634475228
$ws.Cells.Item(2,16).Value = 0.001903591634475228
$ws.Cells.Item(2,17).Value = 0.002171158743555556
$ws.Cells.Item(2,18).Value = 0.019540428692
$ws.Cells.Item(2,19).Value = 0.000004592115266292263
$ws.Cells.Item(2,20).Value = 0.000004592115266292261

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Reln"
$ws.Cells.Item(3,3).Value = "Vldlr"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 0.03565433333333334
$ws.Cells.Item(3,8).Value = 0.106963
$ws.Cells.Item(3,9).Value = 0.002412342638581826
$ws.Cells.Item(3,10).Value = 0.002412342638581825
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 28.046323
$ws.Cells.Item(3,14).Value = 84.138969
$ws.Cells.Item(3,15).Value = 0.8767392739472014
$ws.Cells.Item(3,16).Value = 0.8767392739472013
$ws.Cells.Item(3,17).Value = 0.9999729490163335
$ws.Cells.Item(3,18).Value = 8.999756541147001
$ws.Cells.Item(3,19).Value = 0.002114995533462106
$ws.Cells.Item(3,20).Value = 0.002114995533462105

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Reln"
$ws.Cells.Item(4,3).Value = "Vldlr"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 0.03565433333333334
$ws.Cells.Item(4,8).Value = 0.106963
$ws.Cells.Item(4,9).Value = 0.002412342638581826
$ws.Cells.Item(4,10).Value = 0.002412342638581825
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 3.882136333333333
$ws.Cells.Item(4,14).Value = 11.646409
$ws.Cells.Item(4,15).Value = 0.1213571344183235
$ws.Cells.Item(4,16).Value = 0.1213571344183235
$ws.Cells.Item(4,17).Value = 0.1384149828741111
$ws.Cells.Item(4,18).Value = 1.245734845867
$ws.Cells.Item(4,19).Value = 0.0002927549898534277
$ws.Cells.Item(4,20).Value = 0.0002927549898534277

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Reln"
$ws.Cells.Item(5,3).Value = "Vldlr"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 8.066615333333333
$ws.Cells.Item(5,8).Value = 24.199846
$ws.Cells.Item(5,9).Value = 0.5457805068380079
$ws.Cells.Item(5,10).Value = 0.5457805068380079
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,12).Value = 0.3333333333333333
$ws.Cells.Item(5,13).Value = 0.06089466666666667
$ws.Cells.Item(5,14).Value = 0.182684
$ws.Cells.Item(5,15).Value = 0.001903591634475228
$ws.Cells.Item(5,16).Value = 0.001903591634475228
$ws.Cells.Item(5,17).Value = 0.4912138518515556
$ws.Cells.Item(5,18).Value = 4.420924666664001
$ws.Cells.Item(5,19).Value = 0.001038943207076482
$ws.Cells.Item(5,20).Value = 0.001038943207076482

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Reln"
$ws.Cells.Item(6,3).Value = "Vldlr"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 8.066615333333333
$ws.Cells.Item(6,8).Value = 24.199846
$ws.Cells.Item(6,9).Value = 0.5457805068380079
$ws.Cells.Item(6,10).Value = 0.5457805068380079
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 28.046323
$ws.Cells.Item(6,14).Value = 84.138969
$ws.Cells.Item(6,15).Value = 0.8767392739472014
$ws.Cells.Item(6,16).Value = 0.8767392739472013
$ws.Cells.Item(6,17).Value = 226.2388991554193
$ws.Cells.Item(6,18).Value = 2036.150092398774
$ws.Cells.Item(6,19).Value = 0.4785072052996907
$ws.Cells.Item(6,20).Value = 0.4785072052996906

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Reln"
$ws.Cells.Item(7,3).Value = "Vldlr"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 8.066615333333333
$ws.Cells.Item(7,8).Value = 24.199846
$ws.Cells.Item(7,9).Value = 0.5457805068380079
$ws.Cells.Item(7,10).Value = 0.5457805068380079
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 3.882136333333333
$ws.Cells.Item(7,14).Value = 11.646409
$ws.Cells.Item(7,15).Value = 0.1213571344183235
$ws.Cells.Item(7,16).Value = 0.1213571344183235
$ws.Cells.Item(7,17).Value = 31.31570047255711
$ws.Cells.Item(7,18).Value = 281.841304253014
$ws.Cells.Item(7,19).Value = 0.06623435833124083
$ws.Cells.Item(7,20).Value = 0.06623435833124083

# Row 8
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Reln"
$ws.Cells.Item(8,3).Value = "Vldlr"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 6.677692666666666
$ws.Cells.Item(8,8).Value = 20.033078
$ws.Cells.Item(8,9).Value = 0.4518071505234102
$ws.Cells.Item(8,10).Value = 0.4518071505234102
$ws.Cells.Item(8,11).Value = 1
$ws.Cells.Item(8,12).Value = 0.3333333333333333
$ws.Cells.Item(8,13).Value = 0.06089466666666667
$ws.Cells.Item(8,14).Value = 0.182684
$ws.Cells.Item(8,15).Value = 0.001903591634475228
$ws.Cells.Item(8,16).Value = 0.001903591634475228
$ws.Cells.Item(8,17).Value = 0.4066358690391111
$ws.Cells.Item(8,18).Value = 3.659722821352
$ws.Cells.Item(8,19).Value = 0.0008600563121324538
$ws.Cells.Item(8,20).Value = 0.0008600563121324537

# Row 9
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Reln"
$ws.Cells.Item(9,3).Value = "Vldlr"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 6.677692666666666
$ws.Cells.Item(9,8).Value = 20.033078
$ws.Cells.Item(9,9).Value = 0.4518071505234102
$ws.Cells.Item(9,10).Value = 0.4518071505234102
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 28.046323
$ws.Cells.Item(9,14).Value = 84.138969
$ws.Cells.Item(9,15).Value = 0.8767392739472014
$ws.Cells.Item(9,16).Value = 0.8767392739472013
$ws.Cells.Item(9,17).Value = 187.2847254240647
$ws.Cells.Item(9,18).Value = 1685.562528816582
$ws.Cells.Item(9,19).Value = 0.3961170731140486
$ws.Cells.Item(9,20).Value = 0.3961170731140485

# Row 10
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Reln"
$ws.Cells.Item(10,3).Value = "Vldlr"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 6.677692666666666
$ws.Cells.Item(10,8).Value = 20.033078
$ws.Cells.Item(10,9).Value = 0.4518071505234102
$ws.Cells.Item(10,10).Value = 0.4518071505234102
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 3.882136333333333
$ws.Cells.Item(10,14).Value = 11.646409
$ws.Cells.Item(10,15).Value = 0.1213571344183235
$ws.Cells.Item(10,16).Value = 0.1213571344183235
$ws.Cells.Item(10,17).Value = 25.92371332410022
$ws.Cells.Item(10,18).Value = 233.313419916902
$ws.Cells.Item(10,19).Value = 0.0548300210972292
$ws.Cells.Item(10,20).Value = 0.0548300210972292
